# Refitting NCDEs to individual patients (for manuscript figure)
#
# Adds a new "Label" column (H) that flags each row as Control (0) or
# MDD (1), and updates the refit Prediction/Error/Cross Entropy Loss
# values in columns D:F to the new fitted numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Range("H1").Value = "Label"

# --- Block 1 (Iterations = 100), rows 2-11 ---

# Row 2: Control 39
$ws.Range("D2").Value = 0.9982129799411923
$ws.Range("E2").Value = 0.9982129799411923
$ws.Range("H2").Value = 0

# Row 3: Control 17
$ws.Range("D3").Value = 0.3761847003727704
$ws.Range("E3").Value = 0.3761847003727704
$ws.Range("H3").Value = 0

# Row 4: Control 23
$ws.Range("D4").Value = 0.4779530881917586
$ws.Range("E4").Value = 0.4779530881917586
$ws.Range("H4").Value = 0

# Row 5: Control 27
$ws.Range("D5").Value = 0.7581426521372369
$ws.Range("E5").Value = 0.7581426521372369
$ws.Range("H5").Value = 0

# Row 6: Control 8
$ws.Range("D6").Value = 0.561756656115992
$ws.Range("E6").Value = 0.561756656115992
$ws.Range("H6").Value = 0

# Row 7: MDD 4
$ws.Range("D7").Value = 0.381111600165425
$ws.Range("E7").Value = 0.618888399834575
$ws.Range("H7").Value = 1

# Row 8: MDD 32
$ws.Range("D8").Value = 0.5363363828931194
$ws.Range("E8").Value = 0.4636636171068806
$ws.Range("H8").Value = 1

# Row 9: MDD 48
$ws.Range("D9").Value = 0.4433495877600963
$ws.Range("E9").Value = 0.5566504122399037
$ws.Range("H9").Value = 1

# Row 10: MDD 44
$ws.Range("D10").Value = 0.3690456284863603
$ws.Range("E10").Value = 0.6309543715136396
$ws.Range("H10").Value = 1

# Row 11: MDD 31
$ws.Range("D11").Value = 0.373601980774575
$ws.Range("E11").Value = 0.626398019225425
$ws.Range("F11").Value = 1.407594442367554
$ws.Range("H11").Value = 1

# --- Block 2 (Iterations = 200), rows 12-21 ---
# D/E/F values unchanged in this block per diff; only the Label column is new.

$ws.Range("H12").Value = 0  # Control 39
$ws.Range("H13").Value = 0  # Control 17
$ws.Range("H14").Value = 0  # Control 23
$ws.Range("H15").Value = 0  # Control 27
$ws.Range("H16").Value = 0  # Control 8
$ws.Range("H17").Value = 1  # MDD 4
$ws.Range("H18").Value = 1  # MDD 32
$ws.Range("H19").Value = 1  # MDD 48
$ws.Range("H20").Value = 1  # MDD 44
$ws.Range("H21").Value = 1  # MDD 31
